# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled update).
# Column D ("Price") and E ("Volume(1h)") are plain text in this sheet (not
# numbers), and rows 43/44 swap their Coin/Link/Price/Volume between Maker and
# EnergySwap. Values that look like numbers (e.g. "154.74") are written with a
# leading apostrophe so Excel keeps them as literal text (matching the source's
# exact digits/trailing zeros, e.g. "6.00") instead of auto-converting them to
# numbers; ClearFormats() afterwards drops the transient quote-prefix formatting
# so the cells end up with their normal (default) style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.125.44"
$ws.Range("E2").Value = "  +0.56%  "
$ws.Range("D3").Value = "2.264.69"
$ws.Range("E3").Value = "  -0.11%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'154.74"
$ws.Range("E5").Value = "  +15,353.58%  "
$ws.Range("D6").Value = "'305.91"
$ws.Range("E6").Value = "  +1.34%  "
$ws.Range("D7").Value = "'94.10"
$ws.Range("E7").Value = "  +2.21%  "
$ws.Range("D8").Value = "'0.531"
$ws.Range("E8").Value = "  -0.16%  "
$ws.Range("E9").Value = "  -0.06%  "
$ws.Range("E10").Value = "  +0.18%  "
$ws.Range("D11").Value = "'33.15"
$ws.Range("E11").Value = "  +2.72%  "
$ws.Range("D12").Value = "'0.0802"
$ws.Range("E12").Value = "  +0.42%  "
$ws.Range("E13").Value = "  -2.05%  "
$ws.Range("D14").Value = "'6.67"
$ws.Range("E14").Value = "  +0.14%  "
$ws.Range("D15").Value = "2.615.21"
$ws.Range("E15").Value = "  -0.04%  "
$ws.Range("E16").Value = "  +1.47%  "
$ws.Range("D17").Value = "2.264.11"
$ws.Range("E18").Value = "  +3.61%  "
$ws.Range("D19").Value = "41.985.64"
$ws.Range("E19").Value = "  +0.39%  "
$ws.Range("D20").Value = "'12.68"
$ws.Range("E20").Value = "  +4.65%  "
$ws.Range("E21").Value = "  +1.15%  "
$ws.Range("D22").Value = "'6.00"
$ws.Range("E22").Value = "  +1.01%  "
$ws.Range("D23").Value = "'68.12"
$ws.Range("E23").Value = "  +1.68%  "
$ws.Range("D24").Value = "'244.11"
$ws.Range("E24").Value = "  +0.98%  "
$ws.Range("E25").Value = "  +2.06%  "
$ws.Range("E26").Value = "  +2.12%  "
$ws.Range("D27").Value = "'0.999"
$ws.Range("E27").Value = "  -0.11%  "
$ws.Range("E28").Value = "  +0.27%  "
$ws.Range("D29").Value = "'9.70"
$ws.Range("E29").Value = "  +0.58%  "
$ws.Range("E30").Value = "  +0.30%  "
$ws.Range("D31").Value = "'34.99"
$ws.Range("E31").Value = "  +3.52%  "
$ws.Range("D32").Value = "'159.96"
$ws.Range("E32").Value = "  +0.41%  "
$ws.Range("E33").Value = "  +3.45%  "
$ws.Range("D34").Value = "'1.00"
$ws.Range("E34").Value = "  +0.07%  "
$ws.Range("E35").Value = "  -0.18%  "
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("E37").Value = "  +2.84%  "
$ws.Range("E38").Value = "  -0.83%  "
$ws.Range("E39").Value = "  +1.02%  "
$ws.Range("E40").Value = "  +0.65%  "
$ws.Range("E41").Value = "  -1.26%  "
$ws.Range("D42").Value = "'4.08"
$ws.Range("E42").Value = "  +4.27%  "
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").Value = "'19.63"
$ws.Range("E43").Value = "  -0.41%  "
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "2.006.06"
$ws.Range("E44").Value = "  -3.07%  "
$ws.Range("D45").Value = "'2.26"
$ws.Range("E45").Value = "  +11.05%  "
$ws.Range("E46").Value = "  +1.25%  "
$ws.Range("D47").Value = "'10.19"
$ws.Range("E47").Value = "  -0.11%  "
$ws.Range("D48").Value = "'2.92"
$ws.Range("E48").Value = "  +0.20%  "
$ws.Range("E49").Value = "  +3.68%  "
$ws.Range("D50").Value = "'73.02"
$ws.Range("E50").Value = "  +3.06%  "
$ws.Range("E51").Value = "  +0.30%  "

# Cells above written with a leading apostrophe (text-forced) -- strip the
# resulting quote-prefix number format so no stray style is left behind.
$textForcedCells = @(
    "D5", "D6", "D7", "D8", "D11", "D12",
    "D14", "D20", "D22", "D23", "D24", "D27",
    "D29", "D31", "D32", "D34", "D42", "D43",
    "D45", "D47", "D48", "D50"
)
foreach ($ref in $textForcedCells) {
    $ws.Range($ref).ClearFormats()
}
